$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 20220529
$ws.Range("B2").Value = "DW"
$ws.Range("C2").Value = "Push to github"
$ws.Range("D2").Value = "NA"

$ws.Columns.Item(1).ColumnWidth = 9.375
$ws.Columns.Item(3).ColumnWidth = 16

$ws.Range("C2").Select()
